# The deck was re-saved roughly a week later; PowerPoint recached the
# "datetimeFigureOut" date field that lives on the Date placeholder of
# every slide layout, bumping its displayed text from 5/9/18 to
# 5/16/2018. Walk every custom layout on the slide master and update
# the Date Placeholder shape's text accordingly.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$newDate = "5/16/2018"

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
